$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 20.875
$ws.Range("I8").Value = 20.875
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 62.625
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 76.375
$ws.Range("H129").Value = 7554.8184
$ws.Range("I129").Value = 386.33334
$ws.Range("J129").Value = 10243
$ws.Range("K129").Value = 1159.00002
$ws.Range("L129").Value = 30729
$ws.Range("M129").Value = 3840.99998
$ws.Range("N129").Value = -40729
$ws.Range("H138").Value = 3139.4177
$ws.Range("J138").Value = 3674.93
$ws.Range("L138").Value = 11024.79
$ws.Range("N138").Value = -21304.79
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 3703.6667
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 4555.5
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 4555.5
$ws.Range("M13").Value = -1856
$ws.Range("N13").Value = -4843.5
$ws.Range("H32").Value = 6345.05
$ws.Range("I32").Value = 4572.7227
$ws.Range("K32").Value = 4572.7227
$ws.Range("M32").Value = -4285.7227
$ws.Range("H45").Value = 2411.2
$ws.Range("I45").Value = 2658.8572
$ws.Range("J45").Value = 1833.3334
$ws.Range("K45").Value = 2658.8572
$ws.Range("L45").Value = 1833.3334
$ws.Range("M45").Value = -2281.8572
$ws.Range("N45").Value = -2587.3334
$ws.Range("H74").Value = 530168.7
$ws.Range("J74").Value = 1430002
$ws.Range("L74").Value = 1430002
$ws.Range("N74").Value = -1431750
$ws.Range("H77").Value = 530168.7
$ws.Range("J77").Value = 1430002
$ws.Range("L77").Value = 7150010
$ws.Range("N77").Value = -7158746
$ws.Range("H139").Value = 37753.75
$ws.Range("J139").Value = 37753.75
$ws.Range("L139").Value = 37753.75
$ws.Range("N139").Value = -48033.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("K10").Value = 2000
$ws.Range("M10").Value = -1861
$ws.Range("H31").Value = 1602.8939
$ws.Range("I31").Value = 1086.234
$ws.Range("J31").Value = 2880.9473
$ws.Range("K31").Value = 1086.234
$ws.Range("L31").Value = 2880.9473
$ws.Range("M31").Value = -791.2339999999999
$ws.Range("N31").Value = -3470.9473
$ws.Range("H34").Value = 1602.8939
$ws.Range("I34").Value = 1086.234
$ws.Range("J34").Value = 2880.9473
$ws.Range("K34").Value = 1086.234
$ws.Range("L34").Value = 2880.9473
$ws.Range("M34").Value = -884.2339999999999
$ws.Range("N34").Value = -3284.9473
$ws.Range("H140").Value = 42633.332
$ws.Range("J140").Value = 42633.332
$ws.Range("L140").Value = 42633.332
$ws.Range("N140").Value = -52993.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 550001
$ws.Range("I9").Value = 550001
$ws.Range("K9").Value = 1650003
$ws.Range("M9").Value = -1649779
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("H122").Value = 1091.55
$ws.Range("I122").Value = 554.7143
$ws.Range("J122").Value = 1380.6154
$ws.Range("K122").Value = 4992.428699999999
$ws.Range("L122").Value = 12425.5386
$ws.Range("M122").Value = -2542.428699999999
$ws.Range("N122").Value = -17325.5386
$ws.Range("H131").Value = 2084197.8
$ws.Range("I131").Value = 6667074.5
$ws.Range("J131").Value = 1071.8182
$ws.Range("K131").Value = 20001223.5
$ws.Range("L131").Value = 3215.4546
$ws.Range("M131").Value = -19996183.5
$ws.Range("N131").Value = -13295.4546
$ws.Range("H132").Value = 1601.8182
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1601.8182
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14416.3638
$ws.Range("N132").Value = -19476.3638
$ws.Range("H136").Value = 1083.4
$ws.Range("I136").Value = 966.9
$ws.Range("J136").Value = 1199.9
$ws.Range("K136").Value = 2900.7
$ws.Range("L136").Value = 3599.7
$ws.Range("M136").Value = 2199.3
$ws.Range("N136").Value = -13799.7
$ws.Range("H138").Value = 1124.9375
$ws.Range("J138").Value = 991.5833
$ws.Range("L138").Value = 2974.7499
$ws.Range("N138").Value = -13254.7499
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 11750326
$ws.Range("I3").Value = 11750326
$ws.Range("K3").Value = 11750326
$ws.Range("M3").Value = -11750210
$ws.Range("H12").Value = 30003
$ws.Range("I12").Value = 30003
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 30003
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -29863
$ws.Range("H136").Value = 25775.162
$ws.Range("J136").Value = 25775.162
$ws.Range("L136").Value = 77325.486
$ws.Range("N136").Value = -82425.486
$ws.Range("N12").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 1250
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -328
$ws.Range("N18").Value = -2344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1900
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -187
$ws.Range("N7").Value = -3726
$ws.Range("H15").Value = 7007
$ws.Range("J15").Value = 7007
$ws.Range("L15").Value = 7007
$ws.Range("N15").Value = -7583
$ws.Range("H81").Value = 45456332
$ws.Range("I81").Value = 66668016
$ws.Range("J81").Value = 2728.2856
$ws.Range("K81").Value = 133336032
$ws.Range("L81").Value = 5456.5712
$ws.Range("M81").Value = -133334971
$ws.Range("N81").Value = -7578.5712
$ws.Range("H84").Value = 45456332
$ws.Range("I84").Value = 66668016
$ws.Range("J84").Value = 2728.2856
$ws.Range("K84").Value = 666680160
$ws.Range("L84").Value = 27282.856
$ws.Range("M84").Value = -666674856
$ws.Range("N84").Value = -37890.856
